$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H3").Value = 720
$ws.Range("H4").Value = 720

$ws.Range("H4").Select()
